# Applies the "Finish User login and register" commit:
#   1. Inserts three new bullet paragraphs before the existing
#      "2-way binding:" paragraph (*ngIf, (click), [hidden]).
#   2. Appends a new "Observables" / rxjs block (plus leading/trailing
#      blank paragraphs) after the final "Making HTTP requests" paragraph.
#   3. Refreshes a couple of stale numbering.xml template ids that Word
#      regenerates once the newly-touched list levels are re-rendered.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------

function New-PkgXml($bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function New-ParaXml($pPrXml, $runsXml) {
    $body = "<w:p>"
    if ($pPrXml) { $body += "<w:pPr>$pPrXml</w:pPr>" }
    $body += $runsXml
    $body += "</w:p>"
    return New-PkgXml $body
}

# Inserts a brand-new paragraph (given its pPr/runs OOXML fragments)
# immediately before $targetPara, leaving $targetPara untouched.
function Insert-ParagraphBefore($targetPara, $pPrXml, $runsXml) {
    $targetPara.Range.InsertParagraphBefore()
    # After InsertParagraphBefore, $targetPara now refers to the freshly
    # created (still empty / inherited-formatting) paragraph; the
    # original content shifted one slot later.
    $newRange = $targetPara.Range
    $newRange.InsertXML((New-ParaXml $pPrXml $runsXml))
}

$runPr26 = '<w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr>'

function Run26($text, $preserve) {
    if ($preserve) {
        return '<w:r>' + $runPr26 + '<w:t xml:space="preserve">' + $text + '</w:t></w:r>'
    }
    return '<w:r>' + $runPr26 + '<w:t>' + $text + '</w:t></w:r>'
}

# ---------------------------------------------------------------------
# 1) Three new paragraphs before "2-way binding:"
# ---------------------------------------------------------------------

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd("`r") -eq "2-way binding:") {
        $target = $d.Paragraphs($i)
        break
    }
}

$listPPr = '<w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr>'

# Insert in reverse order: each insertion lands directly above $target,
# so inserting the last-wanted paragraph first yields the correct order.
Insert-ParagraphBefore $target $listPPr (Run26 "[hidden]: conditional display" $false)
Insert-ParagraphBefore $target $listPPr (Run26 "(click): Angular onclick attribute" $false)
Insert-ParagraphBefore $target $listPPr ((Run26 "*ngIf: conditional " $true) + (Run26 "appear" $false))

# ---------------------------------------------------------------------
# 2) New "Observables" block appended after "Making HTTP requests"
# ---------------------------------------------------------------------

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endRange = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$endRange.InsertParagraphAfter()
$stopper = $d.Paragraphs($d.Paragraphs.Count)

$blankPPr = '<w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr>'

$observablesPPr = '<w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:ind w:left="1080"/><w:rPr><w:b/><w:bCs/><w:color w:val="00B050"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr>'
$observablesRun = '<w:r><w:rPr><w:b/><w:bCs/><w:color w:val="00B050"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>Observables</w:t></w:r>'

$boldListPPr0 = '<w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:b/><w:bCs/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr>'
$listPPr3 = '<w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="2"/></w:numPr><w:ind w:left="1440"/><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr>'
$listPPr1 = '<w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:rPr><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr>'

# Build the full run of new paragraphs (pPr, runs) in document order.
$newParas = @(
    @{ pPr = $blankPPr;      runs = "" },
    @{ pPr = $observablesPPr; runs = $observablesRun },
    @{ pPr = $boldListPPr0;  runs = (Run26 "Lazy collection of multiple values over time" $false) },
    @{ pPr = $listPPr3;      runs = (Run26 "Only subscribers can receive the update" $false) },
    @{ pPr = $listPPr3;      runs = (Run26 "When subscribe: " $true) },
    @{ pPr = $listPPr1;      runs = (Run26 "What to do next with data" $false) },
    @{ pPr = $listPPr1;      runs = (Run26 "What to do when error" $false) },
    @{ pPr = $listPPr1;      runs = (Run26 "What to do when success (optional)" $false) },
    @{ pPr = $boldListPPr0;  runs = ((Run26 "W" $false) + (Run26 "e can send it to JavaScript promise" $false)) },
    @{ pPr = $boldListPPr0;  runs = (Run26 "Pipe(): rxjs extension" $false) },
    @{ pPr = $boldListPPr0;  runs = (Run26 "Parent-to-child and Child-to-parent passing data" $false) },
    @{ pPr = $blankPPr;      runs = "" }
)

# Insert in reverse order before the stopper paragraph so the final
# order matches $newParas top-to-bottom; the stopper itself is left
# untouched (it already has the desired "trailing blank paragraph" shape
# only after we've inserted everything above it -- but since the last
# entry in $newParas *is* that blank paragraph, we insert it too and
# simply remove the now-redundant stopper afterwards).
for ($i = $newParas.Count - 1; $i -ge 0; $i--) {
    $entry = $newParas[$i]
    Insert-ParagraphBefore $stopper $entry.pPr $entry.runs
}

# The stopper paragraph is now a leftover empty paragraph (with
# inherited formatting) sitting after our real trailing blank
# paragraph; delete it so the body ends exactly like the target.
$stopper.Range.Delete() | Out-Null
